# CryCompanywiseStockReport_1.xlsx edit
#
# The source data rows got reshuffled: for a handful of adjacent-row groups
# (each group sharing the same item name in column C), the Balance Qty (B),
# Rate (D), Value (E), Qty (F) and Value (G) figures were rotated by one
# position among the rows in the group (a simple swap for 2-row groups, a
# 3-way rotation for the one 3-row group). Column A (Sr No), C (item name)
# and H:M stay put - only B,D,E,F,G move between rows.
#
# This script re-applies that rotation using the live COM object model:
# for each group we snapshot the B/D/E/F/G values of every row, then write
# row[i] <- old row[i-1] (wrapping around), which reproduces the target
# state described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 4, 5, 6, 7)   # B, D, E, F, G

$groups = @(
    @(149, 150),
    @(264, 265),
    @(313, 314),
    @(316, 317, 318),
    @(346, 347),
    @(351, 352),
    @(355, 356),
    @(372, 373),
    @(389, 390),
    @(419, 420),
    @(421, 422),
    @(457, 458),
    @(579, 580),
    @(583, 584),
    @(590, 591),
    @(593, 594),
    @(599, 600),
    @(601, 602),
    @(720, 721),
    @(872, 873)
)

foreach ($group in $groups) {
    $n = $group.Length

    # Snapshot current values of the columns we care about, per row in the group.
    $snapshot = @()
    for ($i = 0; $i -lt $n; $i++) {
        $row = $group[$i]
        $vals = @()
        foreach ($c in $cols) {
            $vals += $ws.Cells.Item($row, $c).Value2
        }
        $snapshot += ,$vals
    }

    # Write row[i] <- snapshot of row[i-1] (previous row in the group, wrapping
    # around for index 0), which is the rotation the diff encodes.
    for ($i = 0; $i -lt $n; $i++) {
        $row = $group[$i]
        $srcIdx = ($i - 1 + $n) % $n
        $srcVals = $snapshot[$srcIdx]
        for ($k = 0; $k -lt $cols.Length; $k++) {
            $ws.Cells.Item($row, $cols[$k]).Value2 = $srcVals[$k]
        }
    }
}
